# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 34 (pushing all subsequent
# rows down by one) and populate it with the latest survey data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(34).Insert()

$ws.Cells.Item(34, 1).Value = 10
$ws.Cells.Item(34, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(34, 3).Value = "La Araucanía"
$ws.Cells.Item(34, 4).Value = 45281
$ws.Cells.Item(34, 5).Value = 9
$ws.Cells.Item(34, 6).Value = 300000000
$ws.Cells.Item(34, 7).Value = "Espárragos"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 500
$ws.Cells.Item(34, 11).Value = 1800
$ws.Cells.Item(34, 12).Value = 1800
$ws.Cells.Item(34, 13).Value = 1800
$ws.Cells.Item(34, 14).Value = "$/kilo"
$ws.Cells.Item(34, 15).Value = "Región del Maule"
$ws.Cells.Item(34, 16).Value = 1800
$ws.Cells.Item(34, 17).Value = 1
$ws.Cells.Item(34, 18).Value = "Hortaliza"
